$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the header values in B1 and C1 ("V+" and "V-")
$b1 = $ws.Range("B1").Value()
$c1 = $ws.Range("C1").Value()
$ws.Range("B1").Value = $c1
$ws.Range("C1").Value = $b1

# Update the view/selection state (scroll so column B is the left-most
# visible column, with C1 as the active/selected cell)
$ws.Range("C1").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
